# "language login sign up"
# Remove the old demo rows (hello / welcome-name), rename the vi_VN locale
# code to vi-VN, and append the new Login / Sign up / Sign in translation
# strings used by the Experience / Login / Sign up screens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the two demo rows ("hello" and "welcome-name") - everything below
#    shifts up by two rows automatically.
$ws.Rows("2:3").Delete()

# 2) The vi_VN locale code becomes vi-VN (row 1, column C).
$ws.Cells.Item(1, 3).Value = "vi-VN"
$ws.Rows(1).RowHeight = 12.75

# 3) Append the new rows for the Experience / Login / Sign up screens.
$newRows = @(
  ,@("Experience", "Experience Awesome chat", "Trải nghiệm Awesome chat")
  ,@("Login", "Login", "Đăng nhập")
  ,@("DoNotHaveAnAccount", "Don't have an account?", "Chưa có tài khoản?")
  ,@("SignUpNow", "Sign up now", "Đăng ký ngay")
  ,@("ForgotYourPassword", "Forgot your password?", "Quên mật khẩu?")
  ,@("SignUp", "Sign up", "Đăng ký")
  ,@("IagreeToThe", "I agree to the", "Tôi đồng ý với")
  ,@("Policies", "policies", "chính sách")
  ,@("And", "and", "và")
  ,@("Terms", "terms", "điều khoản")
  ,@("AlreadyHaveAnAccount", "Already have an account?", "Đã có tài khoản?")
  ,@("SignInNow", "Sign in now", "Đăng nhập ngay")
  ,@("EnterYourPassword", "Enter your password", "Nhập mật khẩu của bạn")
  ,@("EnterYourEmail", "Enter your email", "Nhập email của bạn")
  ,@("Password", "Password", "Mật khẩu")
)

$startRow = 25
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
  $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
  $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
  $ws.Rows($r).RowHeight = 15.75
}

$lastRow = $startRow + $newRows.Count - 1

# 4) Column widths grew slightly for columns A and C.
$ws.Columns("A").ColumnWidth = 29.75
$ws.Columns("C").ColumnWidth = 33.92

# 5) Page setup: portrait, paper size 9 (A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# 6) Selection moves to the last edited cell (B<lastRow>).
$ws.Cells.Item($lastRow, 2).Select()
